$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 4755.4287
$ws.Range("I43").Value = 1999
$ws.Range("J43").Value = 5214.8335
$ws.Range("K43").Value = 1999
$ws.Range("L43").Value = 5214.8335
$ws.Range("M43").Value = -1930
$ws.Range("N43").Value = -5352.8335

$ws.Range("H62").Value = 6828.625
$ws.Range("J62").Value = 8993.75
$ws.Range("L62").Value = 8993.75
$ws.Range("N62").Value = -10241.75

$ws.Range("H65").Value = 6828.625
$ws.Range("J65").Value = 8993.75
$ws.Range("L65").Value = 44968.75
$ws.Range("N65").Value = -51208.75

$ws.Range("H112").Value = 2737.9092
$ws.Range("J112").Value = 2811.7
$ws.Range("L112").Value = 8435.099999999999
$ws.Range("N112").Value = -10651.1

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H21").Value = 1618.4
$ws.Range("I21").Value = 1015
$ws.Range("J21").Value = 2020.6666
$ws.Range("K21").Value = 1015
$ws.Range("L21").Value = 2020.6666
$ws.Range("M21").Value = -641
$ws.Range("N21").Value = -2768.6666

$ws.Range("H32").Value = 6987.826
$ws.Range("I32").Value = 5272.381
$ws.Range("K32").Value = 5272.381
$ws.Range("M32").Value = -4985.381

$ws.Range("H45").Value = 1661.7142
$ws.Range("J45").Value = 1700
$ws.Range("L45").Value = 1700
$ws.Range("N45").Value = -2454

$ws.Range("H61").Value = 4999
$ws.Range("I61").Value = 4999
$ws.Range("J61").Value = 4999
$ws.Range("K61").Value = 4999
$ws.Range("L61").Value = 4999
$ws.Range("M61").Value = -4787
$ws.Range("N61").Value = -5423

$ws.Range("H63").Value = 5236.125
$ws.Range("J63").Value = 6467.9
$ws.Range("L63").Value = 6467.9
$ws.Range("N63").Value = -7839.9

$ws.Range("H66").Value = 5236.125
$ws.Range("J66").Value = 6467.9
$ws.Range("L66").Value = 32339.5
$ws.Range("N66").Value = -39203.5

$ws.Range("H96").Value = 5172
$ws.Range("J96").Value = 5172
$ws.Range("L96").Value = 5172
$ws.Range("N96").Value = -10664

$ws.Range("H97").Value = 2969.5625
$ws.Range("I97").Value = 1616.0769
$ws.Range("J97").Value = 8834.666999999999
$ws.Range("K97").Value = 1616.0769
$ws.Range("L97").Value = 8834.666999999999
$ws.Range("M97").Value = -1120.0769
$ws.Range("N97").Value = -9826.666999999999

$ws.Range("H106").Value = 22333.334
$ws.Range("J106").Value = 22333.334
$ws.Range("L106").Value = 22333.334
$ws.Range("N106").Value = -24857.334

$ws.Range("H136").Value = 4999
$ws.Range("I136").Value = 4999
$ws.Range("J136").Value = 4999
$ws.Range("K136").Value = 14997
$ws.Range("L136").Value = 14997
$ws.Range("M136").Value = -12447
$ws.Range("N136").Value = -20097

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()

$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()

$ws.Range("H134").Value = 3403.4
$ws.Range("I134").Value = 2860.2856
$ws.Range("K134").Value = 8580.856800000001
$ws.Range("M134").Value = -6045.856800000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2878.3
$ws.Range("I134").Value = 2068.65
$ws.Range("K134").Value = 6205.950000000001
$ws.Range("M134").Value = -3670.950000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 631.8333
$ws.Range("I14").Value = 631.8333
$ws.Range("K14").Value = 1895.4999
$ws.Range("M14").Value = -1722.4999

$ws.Range("H31").Value = 375
$ws.Range("J31").Value = 250
$ws.Range("L31").Value = 750
$ws.Range("N31").Value = -1326

$ws.Range("H68").Value = 1494
$ws.Range("I68").Value = 1240.8334
$ws.Range("J68").Value = 1797.8
$ws.Range("K68").Value = 3722.5002
$ws.Range("L68").Value = 5393.4
$ws.Range("M68").Value = -2911.5002
$ws.Range("N68").Value = -7015.4

$ws.Range("H71").Value = 1494
$ws.Range("I71").Value = 1240.8334
$ws.Range("J71").Value = 1797.8
$ws.Range("K71").Value = 11167.5006
$ws.Range("L71").Value = 16180.2
$ws.Range("M71").Value = -7111.500599999999
$ws.Range("N71").Value = -24292.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 48577.855

$ws.Range("H42").Value = 66355

$ws.Range("H115").Value = 66355

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 4418.1816
$ws.Range("I16").Value = 5178.4443
$ws.Range("K16").Value = 5178.4443
$ws.Range("M16").Value = -5008.4443

$ws.Range("H55").Value = 370.6875
$ws.Range("J55").Value = 286.33334
$ws.Range("L55").Value = 286.33334
$ws.Range("N55").Value = -632.33334

$ws.Range("H68").Value = 3479.25
$ws.Range("I68").Value = 1000.6667
$ws.Range("K68").Value = 1000.6667
$ws.Range("M68").Value = -251.6667

$ws.Range("H71").Value = 3479.25
$ws.Range("I71").Value = 1000.6667
$ws.Range("K71").Value = 5003.3335
$ws.Range("M71").Value = -1259.3335

$ws.Range("H82").Value = 3597.8
$ws.Range("I82").Value = 3580.9167
$ws.Range("J82").Value = 3623.125
$ws.Range("K82").Value = 3580.9167
$ws.Range("L82").Value = 3623.125
$ws.Range("M82").Value = -3219.9167
$ws.Range("N82").Value = -4345.125

$ws.Range("H85").Value = 3597.8
$ws.Range("I85").Value = 3580.9167
$ws.Range("J85").Value = 3623.125
$ws.Range("K85").Value = 3580.9167
$ws.Range("L85").Value = 3623.125
$ws.Range("M85").Value = -2332.9167
$ws.Range("N85").Value = -6119.125

$ws.Range("H100").Value = 1999.6666
$ws.Range("I100").Value = 1499.5
$ws.Range("J100").Value = 3000
$ws.Range("K100").Value = 1499.5
$ws.Range("L100").Value = 3000
$ws.Range("M100").Value = -958.5
$ws.Range("N100").Value = -4082

$ws.Range("H136").Value = 3500
$ws.Range("J136").Value = 4000
$ws.Range("L136").Value = 12000
$ws.Range("N136").Value = -17100

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H25").Value = 2463.5
$ws.Range("J25").Value = 2463.5
$ws.Range("L25").Value = 2463.5
$ws.Range("N25").Value = -3049.5

$ws.Range("H31").Value = 0
$ws.Range("J31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("N31").ClearContents()

$ws.Range("H51").Value = 33000
$ws.Range("I51").Value = 33000
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 33000
$ws.Range("L51").Value = 0
$ws.Range("M51").Value = -32490
$ws.Range("N51").ClearContents()

$ws.Range("H126").Value = 2359.2778
$ws.Range("I126").Value = 1797.1333
$ws.Range("J126").Value = 5170
$ws.Range("K126").Value = 5391.3999
$ws.Range("L126").Value = 15510
$ws.Range("M126").Value = -2921.3999
$ws.Range("N126").Value = -20450
